$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='29.172.98'; E='  -0.20%  '},
    @{Row=3; D='1.829.12'; E='  -0.81%  '},
    @{Row=4; D='0.9984'; E='  -0.35%  '},
    @{Row=5; D='242.41'; E='  -0.83%  '},
    @{Row=6; D='0.6224'; E='  -1.07%  '},
    @{Row=7; D='1.000'; E='  -0.33%  '},
    @{Row=8; D='0.07377'; E='  -2.30%  '},
    @{Row=9; D='0.2916'; E='  -1.28%  '},
    @{Row=10; D='23.17'; E='  -1.10%  '},
    @{Row=11; D='0.07677'; E='  -0.62%  '},
    @{Row=12; D='1.824.54'; E='  -1.06%  '},
    @{Row=13; D='4.954'; E='  -1.68%  '},
    @{Row=14; D='0.6673'; E='  -1.95%  '},
    @{Row=15; D='82.50'; E='  -1.16%  '},
    @{Row=16; D='0.000008989'; E='  -3.45%  '},
    @{Row=17; D='5.856'; E='  -2.32%  '},
    @{Row=18; D='29.120.94'; E='  -0.32%  '},
    @{Row=19; D='2.076.11'; E='  -0.91%  '},
    @{Row=20; D='236.65'; E='  +0.90%  '},
    @{Row=21; D='12.45'; E='  -2.37%  '},
    @{Row=22; D='0.9999'; E='  -0.38%  '},
    @{Row=23; D='7.344'; E='  +2.05%  '},
    @{Row=24; D='0.9998'; E='  -0.40%  '},
    @{Row=25; D='158.21'; E='  -1.62%  '},
    @{Row=26; D='0.1413'; E='  +0.55%  '},
    @{Row=27; D='8.508'; E='  -0.88%  '},
    @{Row=28; D='17.65'; E='  -2.09%  '},
    @{Row=29; D='0.06019'; E='  +7.55%  '},
    @{Row=30; D='1.487'; E='  -0.91%  '},
    @{Row=31; D='4.096'; E='  -2.62%  '},
    @{Row=32; D='4.075'; E='  -2.17%  '},
    @{Row=33; D='1.206'; E='  -0.58%  '},
    @{Row=34; D='1.866'; E='  +0.13%  '},
    @{Row=35; D='0.7302'; E='  -3.17%  '},
    @{Row=36; D='1.141'; E='  -1.05%  '},
    @{Row=37; D='2.606'; E='  -2.16%  '},
    @{Row=38; D='2.836'; E='  +2.23%  '},
    @{Row=39; D='1.223.83'; E='  -1.50%  '},
    @{Row=40; D='0.01751'; E='  -2.62%  '},
    @{Row=41; D='6.287'; E='  -5.44%  '},
    @{Row=42; D='0.9175'; E='  +1.31%  '},
    @{Row=43; D='1.000'; E='  -0.24%  '},
    @{Row=44; D='101.81'; E='  -0.92%  '},
    @{Row=45; D='1.983.00'; E='  -0.65%  '},
    @{Row=46; D='64.94'; E='  -2.90%  '},
    @{Row=47; D='0.5047'; E='  -1.13%  '},
    @{Row=48; B='BabyDogeCoin'; C='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D='0.00000000117'; E='  -2.32%  '},
    @{Row=49; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.4023'; E='  -2.16%  '},
    @{Row=50; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.079'; E='  -0.88%  '},
    @{Row=51; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1132'; E='  +1.59%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($u.Row, 3).Value = $u.C
    }
    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}

Write-Output "Updated $($updates.Count) rows"